# Fix "Res edit" QC report formatting
#
# The report was missing the TAXYR and TOWNSHIP columns. Insert two new
# columns right after PARID (in front of the old "Class" column), shifting
# every other column two positions to the right, and give the two new
# header cells the same header styling (bold white-on-fill) as the rest of
# row 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns at B:C; everything from the old column B onward
# shifts right by two columns.
$ws.Range("B1:C1").EntireColumn.Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftToRight)

# New header cells: give them the same header format (bold font / fill /
# border) as the neighboring header cells before putting in their labels.
$ws.Range("D1").Copy()
$ws.Range("B1:C1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Range("B1").Value = "TAXYR"
$ws.Range("C1").Value = "TOWNSHIP"

# Size the two new columns for their header text.
$ws.Columns.Item(2).ColumnWidth = 10
$ws.Columns.Item(3).ColumnWidth = 12.666666666666666
